# Repair template and make some works on report
#
# The header row (row 1) of the "Template" sheet is restructured:
#  - AngleUab/AngleUbc/AngleUca (G1:I1) are renamed to AngUab/AngUbc/AngUca
#  - three new columns AngIab/AngIbc/AngIca are inserted right after Ia/Ib/Ic
#    (pushing cosPhi_A..S1 three columns to the right)
#  - the trailing summary columns Pabc/Qabc/Sabc/P1/Q1/S1 are collapsed down
#    to a single P/Q/S triplet at the very end of the row
#
# Net effect: 3 columns are inserted (before cosPhi_A) and 3 columns are
# removed (Pabc, Qabc, Sabc), so the overall used range stays A1:AJ13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 columns for the new AngI?? triplet, right before cosPhi_A (M1) ---
$ws.Range("M1:O1").EntireColumn.Insert()

# --- Delete the 3 now-redundant *abc summary columns (Pabc, Qabc, Sabc) ---
# After the insert above, the old layout shifted right by 3:
#   old Y1 (Pabc) -> now AB1, old AC1 (Qabc) -> now AF1, old AG1 (Sabc) -> now AJ1
$ws.Range("AB1").EntireColumn.Delete()
$ws.Range("AE1").EntireColumn.Delete()
$ws.Range("AH1").EntireColumn.Delete()

# --- Rewrite the header row texts that changed ---
$ws.Range("G1").Value = "AngUab"
$ws.Range("H1").Value = "AngUbc"
$ws.Range("I1").Value = "AngUca"

$ws.Range("M1").Value = "AngIab"
$ws.Range("N1").Value = "AngIbc"
$ws.Range("O1").Value = "AngIca"

# P1..R1 now hold cosPhi_A/B/C, S1..U1 hold U1/U2/U0, V1..X1 hold I1/I2/I0,
# Y1..AA1 hold Pa/Pb/Pc, AB1..AD1 hold Qa/Qb/Qc, AE1..AG1 hold Sa/Sb/Sc -- all
# unchanged text, just shifted by the column insert/delete above.

# Rename the trailing P1/Q1/S1 labels to the shorter P/Q/S
$ws.Range("AH1").Value = "P"
$ws.Range("AI1").Value = "Q"
$ws.Range("AJ1").Value = "S"

# --- Move the current selection to AJ2 ---
$ws.Range("AJ2").Select()
